$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1 (values), then copy O1's format (bold/border style) onto them
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# For each data row (2-25): swap I<->K and M<->O, and add P=2, Q=2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $ws.Cells.Item($r, 9).Value2 = $kVal
    $ws.Cells.Item($r, 11).Value2 = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O
    $ws.Cells.Item($r, 13).Value2 = $oVal
    $ws.Cells.Item($r, 15).Value2 = $mVal

    $ws.Cells.Item($r, 16).Value2 = 2  # column P
    $ws.Cells.Item($r, 17).Value2 = 2  # column Q
}
